$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.230.52"
$ws.Range("E2").Value = "  +3.98%  "

$ws.Range("D3").Value = "1.784.95"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.47%  "

$ws.Range("D5").Value = "336.44"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").Value = "0.9959"

$ws.Range("D7").Value = "0.3829"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "0.3441"

$ws.Range("D9").Value = "47.63"
$ws.Range("E9").Value = "  -0.66%  "

$ws.Range("D10").Value = "1.158"
$ws.Range("E10").Value = "  -2.46%  "

$ws.Range("D11").Value = "0.07434"
$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("D12").Value = "23.08"
$ws.Range("E12").Value = "  +6.54%  "

$ws.Range("D13").Value = "0.9968"
$ws.Range("E13").Value = "  -0.42%  "

$ws.Range("D14").Value = "6.420"
$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("D15").Value = "1.783.60"
$ws.Range("E15").Value = "  +0.19%  "

$ws.Range("D16").Value = "7.135"
$ws.Range("E16").Value = "  +0.67%  "

$ws.Range("D17").Value = "0.00001084"
$ws.Range("E17").Value = "  -0.69%  "

$ws.Range("D18").Value = "0.06657"
$ws.Range("E18").Value = "  +0.27%  "

$ws.Range("E19").Value = "  -0.68%  "

$ws.Range("D20").Value = "0.9957"
$ws.Range("E20").Value = "  -0.50%  "

$ws.Range("D21").Value = "17.50"
$ws.Range("E21").Value = "  +0.78%  "

$ws.Range("D22").Value = "6.430"
$ws.Range("E22").Value = "  -1.35%  "

$ws.Range("D23").Value = "28.206.80"
$ws.Range("E23").Value = "  +3.89%  "

$ws.Range("E24").Value = "  -1.27%  "

$ws.Range("D25").Value = "2.383"
$ws.Range("E25").Value = "  +0.12%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "20.91"
$ws.Range("E26").Value = "  -0.92%  "

$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "1.432"
$ws.Range("E27").Value = "  -0.76%  "

$ws.Range("D28").Value = "2.421"
$ws.Range("E28").Value = "  -2.98%  "

$ws.Range("D29").Value = "153.83"
$ws.Range("E29").Value = "  -0.84%  "

$ws.Range("D30").Value = "1.986.79"
$ws.Range("E30").Value = "  +0.27%  "

$ws.Range("D31").Value = "135.14"
$ws.Range("E31").Value = "  +0.86%  "

$ws.Range("D32").Value = "6.168"
$ws.Range("E32").Value = "  +2.17%  "

$ws.Range("D33").Value = "3.951"
$ws.Range("E33").Value = "  -0.75%  "

$ws.Range("D34").Value = "0.08814"
$ws.Range("E34").Value = "  +1.90%  "

$ws.Range("E35").Value = "  -1.44%  "

$ws.Range("D36").Value = "0.02428"
$ws.Range("E36").Value = "  +4.45%  "

$ws.Range("D37").Value = "0.6873"
$ws.Range("E37").Value = "  +0.64%  "

$ws.Range("D38").Value = "5.333"
$ws.Range("E38").Value = "  -0.96%  "

$ws.Range("D39").Value = "0.06346"
$ws.Range("E39").Value = "  +0.83%  "

$ws.Range("D40").Value = "0.2182"
$ws.Range("E40").Value = "  +0.33%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.244"
$ws.Range("E41").Value = "  +0.25%  "

$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "1.509"
$ws.Range("E42").Value = "  -7.22%  "

$ws.Range("D43").Value = "8.342"
$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("D44").Value = "14.31"
$ws.Range("E44").Value = "  +0.64%  "

$ws.Range("D45").Value = "0.9958"
$ws.Range("E45").Value = "  -0.43%  "

$ws.Range("D46").Value = "0.6325"
$ws.Range("E46").Value = "  -1.37%  "

$ws.Range("D47").Value = "3.852"
$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("D48").Value = "132.16"
$ws.Range("E48").Value = "  +0.69%  "

$ws.Range("D49").Value = "2.096"
$ws.Range("E49").Value = "  -1.66%  "

$ws.Range("D50").Value = "0.07458"
$ws.Range("E50").Value = "  +5.02%  "

$ws.Range("D51").Value = "1.207"
$ws.Range("E51").Value = "  +8.56%  "

